# Apply "Natmi following Dr Hou advice" update to LR-pairs_lrc2p/Il1a-Il1r1 sheet.
# Rows 2-5: sending cluster changes from M2 to ECs (Il1a -> Il1r1, targets ECs/FAPs/M2/sCs).
# Rows 6-9 (new): sending cluster M2 (Il1a -> Il1r1), same four target clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1a"
$ws.Range("C2").Value = "Il1r1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09179766666666667
$ws.Range("H2").Value = 0.275393
$ws.Range("I2").Value = 0.01443540132615123
$ws.Range("J2").Value = 0.01443540132615123
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.21972033333333
$ws.Range("N2").Value = 39.659161
$ws.Range("O2").Value = 0.2311669015805739
$ws.Range("P2").Value = 0.2311669015805739
$ws.Range("Q2").Value = 1.213539480585889
$ws.Range("R2").Value = 10.921855325273
$ws.Range("S2").Value = 0.003336986997638488
$ws.Range("T2").Value = 0.003336986997638488

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1a"
$ws.Range("C3").Value = "Il1r1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09179766666666667
$ws.Range("H3").Value = 0.275393
$ws.Range("I3").Value = 0.01443540132615123
$ws.Range("J3").Value = 0.01443540132615123
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 34.595189
$ws.Range("N3").Value = 103.785567
$ws.Range("O3").Value = 0.6049494580123129
$ws.Range("P3").Value = 0.6049494580123129
$ws.Range("Q3").Value = 3.175757628092333
$ws.Range("R3").Value = 28.581818652831
$ws.Range("S3").Value = 0.008732688208445413
$ws.Range("T3").Value = 0.008732688208445413

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1a"
$ws.Range("C4").Value = "Il1r1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09179766666666667
$ws.Range("H4").Value = 0.275393
$ws.Range("I4").Value = 0.01443540132615123
$ws.Range("J4").Value = 0.01443540132615123
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110909
$ws.Range("N4").Value = 0.332727
$ws.Range("O4").Value = 0.001939412426354648
$ws.Range("P4").Value = 0.001939412426354648
$ws.Range("Q4").Value = 0.01018118741233333
$ws.Range("R4").Value = 0.091630686711
$ws.Range("S4").Value = 0.00002799619671135406
$ws.Range("T4").Value = 0.00002799619671135406

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il1a"
$ws.Range("C5").Value = "Il1r1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09179766666666667
$ws.Range("H5").Value = 0.275393
$ws.Range("I5").Value = 0.01443540132615123
$ws.Range("J5").Value = 0.01443540132615123
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.261089666666665
$ws.Range("N5").Value = 27.783269
$ws.Range("O5").Value = 0.1619442279807586
$ws.Range("P5").Value = 0.1619442279807586
$ws.Range("Q5").Value = 0.8501464221907776
$ws.Range("R5").Value = 7.651317799716999
$ws.Range("S5").Value = 0.00233772992335598
$ws.Range("T5").Value = 0.00233772992335598

# Row 6: M2 -> ECs
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Il1a"
$ws.Range("C6").Value = "Il1r1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.267406666666666
$ws.Range("H6").Value = 18.80222
$ws.Range("I6").Value = 0.9855645986738488
$ws.Range("J6").Value = 0.9855645986738487
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.21972033333333
$ws.Range("N6").Value = 39.659161
$ws.Range("O6").Value = 0.2311669015805739
$ws.Range("P6").Value = 0.2311669015805739
$ws.Range("Q6").Value = 82.85336334860222
$ws.Range("R6").Value = 745.6802701374199
$ws.Range("S6").Value = 0.2278299145829354
$ws.Range("T6").Value = 0.2278299145829354

# Row 7: M2 -> FAPs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Il1a"
$ws.Range("C7").Value = "Il1r1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.267406666666666
$ws.Range("H7").Value = 18.80222
$ws.Range("I7").Value = 0.9855645986738488
$ws.Range("J7").Value = 0.9855645986738487
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 34.595189
$ws.Range("N7").Value = 103.785567
$ws.Range("O7").Value = 0.6049494580123129
$ws.Range("P7").Value = 0.6049494580123129
$ws.Range("Q7").Value = 216.8221181731933
$ws.Range("R7").Value = 1951.39906355874
$ws.Range("S7").Value = 0.5962167698038675
$ws.Range("T7").Value = 0.5962167698038675

# Row 8: M2 -> M2
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Il1a"
$ws.Range("C8").Value = "Il1r1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.267406666666666
$ws.Range("H8").Value = 18.80222
$ws.Range("I8").Value = 0.9855645986738488
$ws.Range("J8").Value = 0.9855645986738487
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.110909
$ws.Range("N8").Value = 0.332727
$ws.Range("O8").Value = 0.001939412426354648
$ws.Range("P8").Value = 0.001939412426354648
$ws.Range("Q8").Value = 0.6951118059933332
$ws.Range("R8").Value = 6.256006253939999
$ws.Range("S8").Value = 0.001911416229643294
$ws.Range("T8").Value = 0.001911416229643293

# Row 9: M2 -> sCs
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Il1a"
$ws.Range("C9").Value = "Il1r1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.267406666666666
$ws.Range("H9").Value = 18.80222
$ws.Range("I9").Value = 0.9855645986738488
$ws.Range("J9").Value = 0.9855645986738487
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.261089666666665
$ws.Range("N9").Value = 27.783269
$ws.Range("O9").Value = 0.1619442279807586
$ws.Range("P9").Value = 0.1619442279807586
$ws.Range("Q9").Value = 58.04301511746443
$ws.Range("R9").Value = 522.3871360571799
$ws.Range("S9").Value = 0.1596064980574026
$ws.Range("T9").Value = 0.1596064980574026
